# Commit: "Updated Excel file and save_reqs logic"
#
# Changes:
#   1. Rename the "Steam Turbine" sheet tab to "Steam Turbine " (trailing space).
#   2. Make "Steam Turbine" the active sheet (was "Key"), which also moves the
#      sheet's selection/cursor, updating Steam Turbine's saved selection to E15.
#   3. The previously-active "Key" sheet keeps its own saved selection (B11) but
#      is no longer the active/tabSelected sheet.

$wb = $excel.ActiveWorkbook

$steamTurbine = $wb.Worksheets.Item("Steam Turbine")

# Rename "Steam Turbine" -> "Steam Turbine " (trailing space added)
$steamTurbine.Name = "Steam Turbine "

# Make it the active sheet, and move its selection to E15
$steamTurbine.Activate()
$steamTurbine.Range("E15").Select()
